$d = $word.ActiveDocument

# Change 1: merge ", " + "OLAP Data Models" runs into one run text (no visible text change)
$d.Content.Find.Execute(", OLAP Data Models", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ", OLAP Data Models", 2) | Out-Null

# Change 2: fix "Exercsed" -> "Exercise"
$d.Content.Find.Execute("Exercsed", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Exercise", 2) | Out-Null
